$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric need to be forced to text so Excel
# does not auto-convert them (the workbook stores these columns as text).
$textForceCells = @("D4", "D5", "D6", "D9", "D11", "D12", "D16", "D17", "D19", "D20", "D22", "D23", "D24", "D27", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D41", "D42", "D45", "D47", "D49", "D50", "D51")
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "42.834.98"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.529.81"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "311.42"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "100.91"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "2.921.68"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "2.570.88"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "15.38"
$ws.Range("E16").Value = "  -3.29%  "
$ws.Range("D17").Value = "0.815"
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "42.811.20"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "6.69"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "12.38"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "0.0₃0953"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "69.73"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "243.80"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D27").Value = "25.49"
$ws.Range("E27").Value = "  -5.58%  "
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "10.19"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "38.78"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("D31").Value = "162.19"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "5.80"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +8.82%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "18.37"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  -6.64%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  -6.00%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "4.18"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").Value = "22.11"
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "0.0300"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "1.986.68"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "9.20"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("D48").Value = "2.775.41"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "0.191"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("D50").Value = "79.54"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "72.44"
$ws.Range("E51").Value = "  -2.21%  "

# Restore default (General/Normal) style on the cells we forced to text,
# so no stray cell-level style index is left on them.
foreach ($ref in $textForceCells) {
    $ws.Range($ref).Style = "Normal"
}
